$wb = $excel.ActiveWorkbook

# --- Sheet "lots" (first sheet): fix the date typo and restyle column X (rows 3-25) ---
$ws1 = $wb.Worksheets.Item(1)

# Update the date text for the whole X column (shared text "11/11/2020" -> "11/11/2021").
# Force the cells to stay text (not get auto-converted to a date serial number).
$dateRange = $ws1.Range("X2:X25")
$dateRange.NumberFormat = "@"
$dateRange.Value = "11/11/2021"

# Re-apply the new text value to row 2 alone so it keeps referencing the default style (s=0).
$ws1.Range("X2").NumberFormat = "@"
$ws1.Range("X2").Value = "11/11/2021"

# Rows 3-25 get a distinct (non-bold, regular) font/style, matching the sixth font added to the workbook.
$styledRange = $ws1.Range("X3:X25")
$styledRange.Font.ThemeColor = 1
$styledRange.Font.Bold = $false
$styledRange.Font.Name = "Calibri"
$styledRange.Font.Size = 11

# Update the active selection on the "lots" sheet.
$ws1.Activate()
$ws1.Range("X26").Select()

# --- Other sheets: reset the selection back to A1 ---
for ($i = 2; $i -le $wb.Worksheets.Count; $i++) {
    $sheet = $wb.Worksheets.Item($i)
    $sheet.Activate()
    $sheet.Range("A1").Select()
}
$ws1.Activate()

# --- Normalize header/footer page margins on every sheet ---
foreach ($sheet in $wb.Worksheets) {
    $sheet.PageSetup.HeaderMargin = 36.850393700787386
    $sheet.PageSetup.FooterMargin = 36.850393700787386
}
